# Update cryptos list with latest prices/volumes (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.414.99"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.77%  '

$ws.Range('D3').Value = "'1.796.44"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.99%  '

$ws.Range('D4').Value = "'1.007"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('E5').Value = '  +0.07%  '

$ws.Range('D6').Value = "'307.50"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '

$ws.Range('D7').Value = "'0.4546"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.78%  '

$ws.Range('D8').Value = "'0.3594"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.05%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "'0.07119"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.68%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = "'0.8862"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.03%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = "'0.07823"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.02%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = "'19.51"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.48%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.804.84"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.34%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'5.287"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.05%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'6.333"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.94%  '

$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = "'85.18"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.08%  '

$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = "'1.008"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.12%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = "'0.000008592"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.59%  '

$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = "'1.006"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.08%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = "'14.29"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.14%  '

$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = "'26.437.79"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.76%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'4.991"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.29%  '

$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = "'2.027.61"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.22%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = "'10.56"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.10%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = "'1.984"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.41%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'152.78"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.23%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'17.91"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.86%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = "'2.050"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.03%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'111.99"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.48%  '

$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'4.874"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.55%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = "'0.08659"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.95%  '

$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').Value = "'3.057"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.24%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'4.453"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.11%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'0.7280"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.76%  '

$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = "'2.716"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.36%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'1.112"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.42%  '

$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = "'1.076"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.93%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'0.01946"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.51%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = "'0.05117"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.30%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = "'2.876"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.67%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'0.5176"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.17%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'6.897"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.13%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = "'0.1524"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.67%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'8.016"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.56%  '

$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.4674"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.13%  '

$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = "'1.007"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.12%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'9.852"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.17%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'100.20"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.19%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'1.588"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.57%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.05976"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.99%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'64.23"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.75%  '

